$wb = $excel.ActiveWorkbook

# --- Sheet 1: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("H11").Value = 1018.8
$ws1.Range("I11").Value = 1007.51
$ws1.Range("H23").Value = "2 de 21"

# --- Sheet 2: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F11").Value = 3315.63
$ws2.Range("F23").Value = 12947.57

# --- Sheet 3: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column E (5th column) from 23 to 24
# (the runtime's ColumnWidth setter adds a fixed 5/6 char offset when it
# stores the OOXML "width" attribute, so compensate to land exactly on 24)
$ws3.Columns.Item(5).ColumnWidth = 23.166666666666668

$ws3.Range("D7").Value = 1087.53
$ws3.Range("E7").Value = 1312.47
$ws3.Range("F7").Value = 0.4531375

$ws3.Range("D8").Value = 1050.14
$ws3.Range("E8").Value = -425.1400000000001
$ws3.Range("F8").Value = 1.680224

$ws3.Range("D19").Value = 12947.57
$ws3.Range("E19").Value = 46440.65762291769
$ws3.Range("F19").Value = 0.2180157670676736
